# The author deleted two rows from the apartment-complex list:
#   row 20 -> "뉴타운4차" (id 13924)
#   row 21 -> "뉴타운5차" (id 1995)
# Deleting both rows at once shifts every row below them up by two,
# which is exactly what the target diff shows (rows 22.. become 20..,
# and the previously-last rows 114/115 fall off the bottom). Excel also
# drops the now-unused shared-string entries for the deleted names on
# save, matching the sharedStrings.xml count/uniqueCount change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20:B21").EntireRow.Delete()

# Leave the selection where the author's final screenshot/save shows it.
$ws.Range("G20").Select()
